$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 17.49117535635123
$ws.Cells.Item(2, 4).Value = 10.10690019881205
$ws.Cells.Item(2, 5).Value = 13.61178519483119
$ws.Cells.Item(2, 6).Value = 52.59671278604533
$ws.Cells.Item(2, 7).Value = 3.728657829120948
$ws.Cells.Item(2, 11).Value = 11.80445246554553
$ws.Cells.Item(2, 12).Value = 9.619084943955858
$ws.Cells.Item(3, 2).Value = 17.43832787413744
$ws.Cells.Item(3, 4).Value = 9.964007593499755
$ws.Cells.Item(3, 5).Value = 13.5572864922265
$ws.Cells.Item(3, 6).Value = 51.39187863642607
$ws.Cells.Item(3, 7).Value = 3.732758932429832
$ws.Cells.Item(3, 11).Value = 11.75518689163667
$ws.Cells.Item(3, 12).Value = 9.62190991643773
$ws.Cells.Item(4, 2).Value = 17.41053870755996
$ws.Cells.Item(4, 4).Value = 9.87452319141482
$ws.Cells.Item(4, 5).Value = 13.52288642561479
$ws.Cells.Item(4, 6).Value = 50.64115787048534
$ws.Cells.Item(4, 7).Value = 3.735403407206262
$ws.Cells.Item(4, 11).Value = 11.73182095276525
$ws.Cells.Item(4, 12).Value = 9.62635256973878
$ws.Cells.Item(5, 2).Value = 17.40039431402744
$ws.Cells.Item(5, 4).Value = 9.837636159600482
$ws.Cells.Item(5, 5).Value = 13.50863149308258
$ws.Cells.Item(5, 6).Value = 50.33281286931101
$ws.Cells.Item(5, 7).Value = 3.736512965059339
$ws.Cells.Item(5, 11).Value = 11.72403757076639
$ws.Cells.Item(5, 12).Value = 9.62884328105353
$ws.Cells.Item(6, 2).Value = 17.39878132398184
$ws.Cells.Item(6, 4).Value = 9.831486178619697
$ws.Cells.Item(6, 5).Value = 13.50625015573252
$ws.Cells.Item(6, 6).Value = 50.2814762329448
$ws.Cells.Item(6, 7).Value = 3.73669913766919
$ws.Cells.Item(6, 11).Value = 11.72285029025888
$ws.Cells.Item(6, 12).Value = 9.629297921893322
$ws.Cells.Item(7, 2).Value = 17.41039710942732
$ws.Cells.Item(7, 4).Value = 9.874027402060419
$ws.Cells.Item(7, 5).Value = 13.52269513771313
$ws.Cells.Item(7, 6).Value = 50.63700877177932
$ws.Cells.Item(7, 7).Value = 3.735418241703536
$ws.Cells.Item(7, 11).Value = 11.73170893745495
$ws.Cells.Item(7, 12).Value = 9.626383407095471
$ws.Cells.Item(8, 2).Value = 17.47199186102468
$ws.Cells.Item(8, 4).Value = 10.058004302577
$ws.Cells.Item(8, 5).Value = 13.59318813438514
$ws.Cells.Item(8, 6).Value = 52.1837533594526
$ws.Cells.Item(8, 7).Value = 3.730045739767564
$ws.Cells.Item(8, 11).Value = 11.78604274918342
$ws.Cells.Item(8, 12).Value = 9.619496636207574
$ws.Cells.Item(9, 2).Value = 17.62932375518433
$ws.Cells.Item(9, 4).Value = 10.40414489179929
$ws.Cells.Item(9, 5).Value = 13.72402619775453
$ws.Cells.Item(9, 6).Value = 55.116199173488
$ws.Cells.Item(9, 7).Value = 3.720506956406427
$ws.Cells.Item(9, 11).Value = 11.94669030706364
$ws.Cells.Item(9, 12).Value = 9.627492732499194
$ws.Cells.Item(10, 2).Value = 17.76653685082715
$ws.Cells.Item(10, 4).Value = 10.64853931559253
$ws.Cells.Item(10, 5).Value = 13.81570428845927
$ws.Cells.Item(10, 6).Value = 57.19190766522922
$ws.Cells.Item(10, 7).Value = 3.71409786129467
$ws.Cells.Item(10, 11).Value = 12.09674868346444
$ws.Cells.Item(10, 12).Value = 9.646468410151936
$ws.Cells.Item(11, 2).Value = 17.83348677055669
$ws.Cells.Item(11, 4).Value = 10.75736838639029
$ws.Cells.Item(11, 5).Value = 13.85645835809923
$ws.Cells.Item(11, 6).Value = 58.1157840179516
$ws.Cells.Item(11, 7).Value = 3.711310441889169
$ws.Cells.Item(11, 11).Value = 12.171695180417
$ws.Cells.Item(11, 12).Value = 9.657935217303198
$ws.Cells.Item(12, 2).Value = 17.85947457260501
$ws.Cells.Item(12, 4).Value = 10.79822553236517
$ws.Cells.Item(12, 5).Value = 13.87175506305641
$ws.Cells.Item(12, 6).Value = 58.46245042259623
$ws.Cells.Item(12, 7).Value = 3.710273198958681
$ws.Cells.Item(12, 11).Value = 12.20100917687368
$ws.Cells.Item(12, 12).Value = 9.66268336148914
$ws.Cells.Item(13, 2).Value = 17.85384964853293
$ws.Cells.Item(13, 4).Value = 10.78944221803267
$ws.Cells.Item(13, 5).Value = 13.8684666798444
$ws.Cells.Item(13, 6).Value = 58.38793515524439
$ws.Cells.Item(13, 7).Value = 3.710495776291354
$ws.Cells.Item(13, 11).Value = 12.19465484181467
$ws.Cells.Item(13, 12).Value = 9.661642742625833
$ws.Cells.Item(14, 2).Value = 17.83561216679431
$ws.Cells.Item(14, 4).Value = 10.76073691659699
$ws.Cells.Item(14, 5).Value = 13.85771954816558
$ws.Cells.Item(14, 6).Value = 58.14436971115737
$ws.Cells.Item(14, 7).Value = 3.711224741426748
$ws.Cells.Item(14, 11).Value = 12.17408832292983
$ws.Cells.Item(14, 12).Value = 9.658317726284549
$ws.Cells.Item(15, 2).Value = 17.82452341950395
$ws.Cells.Item(15, 4).Value = 10.74310750216253
$ws.Cells.Item(15, 5).Value = 13.85111893342559
$ws.Cells.Item(15, 6).Value = 57.99475673139505
$ws.Cells.Item(15, 7).Value = 3.711673631969351
$ws.Cells.Item(15, 11).Value = 12.1616114192887
$ws.Cells.Item(15, 12).Value = 9.656333861156702
$ws.Cells.Item(16, 2).Value = 17.76225132938685
$ws.Cells.Item(16, 4).Value = 10.64137841935413
$ws.Cells.Item(16, 5).Value = 13.8130218459925
$ws.Cells.Item(16, 6).Value = 57.13109789555632
$ws.Cells.Item(16, 7).Value = 3.714282592513389
$ws.Cells.Item(16, 11).Value = 12.09198294369604
$ws.Cells.Item(16, 12).Value = 9.645775932479788
$ws.Cells.Item(17, 2).Value = 17.72519872538503
$ws.Cells.Item(17, 4).Value = 10.57835777239147
$ws.Cells.Item(17, 5).Value = 13.78940753468287
$ws.Cells.Item(17, 6).Value = 56.59586486011492
$ws.Cells.Item(17, 7).Value = 3.71591582497939
$ws.Cells.Item(17, 11).Value = 12.05096122222448
$ws.Cells.Item(17, 12).Value = 9.640023941524413
$ws.Cells.Item(18, 2).Value = 17.70431457155871
$ws.Cells.Item(18, 4).Value = 10.54188998321849
$ws.Cells.Item(18, 5).Value = 13.77573567405141
$ws.Cells.Item(18, 6).Value = 56.2861093677222
$ws.Cells.Item(18, 7).Value = 3.716867283176587
$ws.Cells.Item(18, 11).Value = 12.02799761080893
$ws.Cells.Item(18, 12).Value = 9.636982492956518
$ws.Cells.Item(19, 2).Value = 17.69731745965104
$ws.Cells.Item(19, 4).Value = 10.52950534970461
$ws.Cells.Item(19, 5).Value = 13.77109121692666
$ws.Cells.Item(19, 6).Value = 56.18091264217935
$ws.Cells.Item(19, 7).Value = 3.717191506889332
$ws.Cells.Item(19, 11).Value = 12.02033170227461
$ws.Cells.Item(19, 12).Value = 9.635998603326312
$ws.Cells.Item(20, 2).Value = 17.72909890526634
$ws.Cells.Item(20, 4).Value = 10.58508931687792
$ws.Cells.Item(20, 5).Value = 13.79193057940802
$ws.Cells.Item(20, 6).Value = 56.65304026029528
$ws.Cells.Item(20, 7).Value = 3.715740716742397
$ws.Cells.Item(20, 11).Value = 12.05526294125705
$ws.Cells.Item(20, 12).Value = 9.640608633990057
$ws.Cells.Item(21, 2).Value = 17.84095184966656
$ws.Cells.Item(21, 4).Value = 10.76917809701033
$ws.Cells.Item(21, 5).Value = 13.86087992575255
$ws.Cells.Item(21, 6).Value = 58.21599918361758
$ws.Cells.Item(21, 7).Value = 3.711010131282395
$ws.Cells.Item(21, 11).Value = 12.18010411115319
$ws.Cells.Item(21, 12).Value = 9.659283363688928
$ws.Cells.Item(22, 2).Value = 17.9177490598935
$ws.Cells.Item(22, 4).Value = 10.88741983380607
$ws.Cells.Item(22, 5).Value = 13.90515007390252
$ws.Cells.Item(22, 6).Value = 59.21881040704257
$ws.Cells.Item(22, 7).Value = 3.708024981639128
$ws.Cells.Item(22, 11).Value = 12.26712078446017
$ws.Cells.Item(22, 12).Value = 9.67385323955487
$ws.Cells.Item(23, 2).Value = 17.87642850772387
$ws.Cells.Item(23, 4).Value = 10.8245067693667
$ws.Cells.Item(23, 5).Value = 13.88159450876964
$ws.Cells.Item(23, 6).Value = 58.68537929737536
$ws.Cells.Item(23, 7).Value = 3.709608504641921
$ws.Cells.Item(23, 11).Value = 12.22019175540706
$ws.Cells.Item(23, 12).Value = 9.665861319436129
$ws.Cells.Item(24, 2).Value = 17.72733433111584
$ws.Cells.Item(24, 4).Value = 10.58204672062354
$ws.Cells.Item(24, 5).Value = 13.79079020896058
$ws.Cells.Item(24, 6).Value = 56.62719760540767
$ws.Cells.Item(24, 7).Value = 3.715819844259745
$ws.Cells.Item(24, 11).Value = 12.05331619984001
$ws.Cells.Item(24, 12).Value = 9.640343467592002
$ws.Cells.Item(25, 2).Value = 17.58291160200008
$ws.Cells.Item(25, 4).Value = 10.31217608426111
$ws.Cells.Item(25, 5).Value = 13.68941667310752
$ws.Cells.Item(25, 6).Value = 54.33553788071104
$ws.Cells.Item(25, 7).Value = 3.722981633796741
$ws.Cells.Item(25, 11).Value = 11.89752344161442
$ws.Cells.Item(25, 12).Value = 9.623026685037257
